# 2020.05.24 update: add 李泓烨's reading report row, and fill in the
# "他人提出的问题" scores (column E) for 方骏 / 袁佳怡 / 孙毅远 / 闻浩.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- fill in missing 分数（35分） scores in column E ---------------------
$ws.Range("E4").Value = 30
$ws.Range("E5").Value = 35
$ws.Range("E6").Value = 30

# --- insert a new data row for 李泓烨 above the 闻浩 row (row 8) ---------
$ws.Rows.Item(8).Insert()

$ws.Range("A8").Value = "李泓烨"
$ws.Range("B8").Value = 0
$ws.Range("D8").Value = 2
$ws.Range("E8").Value = 35
$ws.Range("F8").Value = "看完第12章"
$ws.Range("G8").Value = 15
$ws.Range("H8").Value = "无"
$ws.Range("I8").Value = 0
$ws.Range("K8").Value = 6
$ws.Range("L8").Formula = "=SUM(C8,E8,G8,I8,K8)"

$ws.Range("A8:L8").HorizontalAlignment = -4108
$ws.Range("A8:L8").VerticalAlignment = -4108

# --- 闻浩's row (now shifted down to row 9) also gets its column E score -
$ws.Range("E9").Value = 30

$ws.Range("L10").Select()
